$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 0
$ws1.Range("F8").Value = 0
$ws1.Range("F10").Value = 475

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1609
$ws4.Range("F6").Value = 23
$ws4.Range("F7").Value = 411
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 63
